# correctif d'un oubli (prise en compte du bon nom et optimisation du code)
#
# 1) Capitalize the student names in column A (rows 2-40): "sloth" -> "Sloth",
#    "skido" -> "Skido", etc.
# 2) Column F was a leftover "date corrected" placeholder ("le 05/16 a 11h53")
#    styled with the (unused/incorrect) green fill. It should really just say
#    "Corrigé", styled the same way the "Corrigé" cells in columns C/D are
#    (tomato fill, right aligned) - row 3 already only uses "Non rendu" values
#    and is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Sloth", "Skido", "Alexis", "Amelie", "Baptiste", "Axelle", "Bastian",
    "Camille", "Charlie", "Emeline", "Emile", "Enzo", "Florian", "Gregoire",
    "Hugo", "Jason", "Jean", "Juan", "Louis", "Lea", "Leo", "Manon",
    "Margauxl", "Margauxq", "Marie", "Mateo", "Max", "Maxence", "Mael",
    "Pierre", "Romain", "Sarah", "Theog", "Theob", "Thibaut", "Vivien",
    "Ambre", "Margauxo", "Mona"
)

for ($i = 0; $i -lt $names.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
}

# Reference cell that already carries the correct "Corrigé" style (tomato
# fill + right alignment) used elsewhere on the sheet (columns C/D).
$corrige = $ws.Range("C2")

foreach ($row in 2..40) {
    if ($row -eq 3) { continue }
    $fcell = $ws.Cells.Item($row, 6)
    $fcell.Value = "Corrigé"
    $fcell.Interior.Color = $corrige.Interior.Color
    $fcell.HorizontalAlignment = $corrige.HorizontalAlignment
}
